# The document has two logos (Pearson logo in the footers, BTEC logo in the
# headers) embedded as inline pictures. Each one is duplicated across the
# "odd" and "even" Header/Footer parts. The edit swaps the cosmetic picture
# "name" (wp:docPr / pic:cNvPr @name) that Word shows in the Selection Pane:
#   - Pearson logo pictures (alt text contains "PearsonLogo"): image1.png -> image2.png
#   - BTEC logo pictures   (alt text contains "BTec_Logo-Orange"): image2.jpg -> image1.jpg
#
# InlineShape has no settable "Name" property (matching real Word's object
# model - only floating Shape objects expose Name), so each picture is
# temporarily converted to a floating Shape, renamed, then converted back to
# an inline picture so the wp:inline layout is preserved.

$d = $word.ActiveDocument

function Rename-LogoPicture($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    for ($i = 1; $i -le 3; $i++) {
        $footer = $section.Footers($i)
        if ($footer.Exists) {
            $shapes = $footer.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $pic = $shapes.Item($j)
                if ($pic.AlternativeText -like "*PearsonLogo*") {
                    Rename-LogoPicture $pic "image2.png"
                }
            }
        }

        $header = $section.Headers($i)
        if ($header.Exists) {
            $shapes = $header.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $pic = $shapes.Item($j)
                if ($pic.AlternativeText -like "*BTec_Logo-Orange*") {
                    Rename-LogoPicture $pic "image1.jpg"
                }
            }
        }
    }
}
